# Apply cryptos list price/volume update (GitHub Actions scraped data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.014.84"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "'1.910.24"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  +5.05%  "

$ws.Range("D6").Value = "'242.07"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.3166"
$ws.Range("E8").Value = "  +2.86%  "

$ws.Range("D9").Value = "'26.38"
$ws.Range("E9").Value = "  +3.07%  "

$ws.Range("D10").Value = "'0.06897"
$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("D11").Value = "'0.08004"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7443"
$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.905.06"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").Value = "'5.190"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "'93.07"
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("D16").Value = "'30.006.16"

$ws.Range("D17").Value = "'13.96"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "'5.869"
$ws.Range("E18").Value = "  -5.36%  "

$ws.Range("D19").Value = "'245.84"
$ws.Range("E19").Value = "  +3.36%  "

$ws.Range("D20").Value = "'0.000007745"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'2.153.56"
$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "'6.837"
$ws.Range("E24").Value = "  -3.32%  "

$ws.Range("D25").Value = "'168.00"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").Value = "'9.230"
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("D27").Value = "'0.1393"
$ws.Range("E27").Value = "  +8.20%  "

$ws.Range("D28").Value = "'18.92"
$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("D29").Value = "'2.034"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").Value = "'1.364"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").Value = "'1.514"
$ws.Range("E31").Value = "  -1.04%  "

$ws.Range("D32").Value = "'4.318"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.081"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.05535"
$ws.Range("E34").Value = "  +3.04%  "

$ws.Range("D35").Value = "'1.257"
$ws.Range("E35").Value = "  -2.45%  "

$ws.Range("D36").Value = "'0.7342"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").Value = "'0.01926"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("D39").Value = "'2.785"
$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("D40").Value = "'6.137"
$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("D41").Value = "'0.4418"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("D42").Value = "'72.50"
$ws.Range("E42").Value = "  -0.50%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "'0.8371"
$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("D45").Value = "'1.879"
$ws.Range("E45").Value = "  -3.83%  "

$ws.Range("D46").Value = "'100.49"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("D47").Value = "'7.557"
$ws.Range("E47").Value = "  -2.04%  "

$ws.Range("D48").Value = "'989.01"
$ws.Range("E48").Value = "  +7.76%  "

$ws.Range("D49").Value = "'2.052.28"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").Value = "'36.26"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("D51").Value = "'1.477"
$ws.Range("E51").Value = "  -0.53%  "
